$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.130.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "'2.636.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'596.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").Value = "'154.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'0.544"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "'2.635.22"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("E10").Value = "  +7.69%  "
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("D12").Value = "'5.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("E14").Value = "  +0.85%  "
$ws.Range("E15").Value = "  +2.40%  "
$ws.Range("D16").Value = "'3.116.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").Value = "'68.035.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").Value = "'2.623.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").Value = "'11.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").Value = "'362.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.28%  "
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("E22").Value = "  +3.04%  "
$ws.Range("D23").Value = "'4.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("E24").Value = "  -1.60%  "
$ws.Range("D25").Value = "'74.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.74%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'9.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.65%  "
$ws.Range("E28").Value = "  +1.00%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("D31").Value = "'559.48"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.07%  "
$ws.Range("D32").Value = "'7.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.31%  "
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("E34").Value = "  +1.07%  "
$ws.Range("E35").Value = "  +1.56%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("E37").Value = "  +3.30%  "
$ws.Range("D38").Value = "'161.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("E39").Value = "  +0.98%  "
$ws.Range("E40").Value = "  +1.33%  "
$ws.Range("E41").Value = "  -0.99%  "
$ws.Range("D42").Value = "'5.29"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.39%  "
$ws.Range("D43").Value = "'0.0₆0341"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.99%  "
$ws.Range("E44").Value = "  +2.47%  "
$ws.Range("E45").Value = "  -1.22%  "
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").Value = "'40.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.72%  "
$ws.Range("D48").Value = "'158.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.03%  "
$ws.Range("E49").Value = "  +1.54%  "
$ws.Range("D50").Value = "'21.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("E51").Value = "  +1.41%  "
